# Apply the "graph2_Vietnam" edit:
#  1. Delete the "United States" sheet.
#  2. On each remaining sheet (Australia, Germany, India):
#       - rename "Wafer" -> "Wafer (excl. polysilicon)"
#       - rename "Cell Cost" -> "Cell Cost (excl. wafer)"
#       - rename "Other material" -> "Other material (e.g. front and back
#         glass, encapsulant and others)"
#       - move that "Other material" row from its old position (row 12,
#         right before "ESG Certification") up to row 6 (right after
#         "Cell Cost"), pushing Overheads/Electricity/Building and
#         facilities/Equipment depreciation/Maintenance/Labour down by one
#         row each
#       - restore the number format on the relocated label cell so it
#         matches the rest of the label column after the cut/insert
#  3. Tidy up the selected cell on each sheet to match the saved view.

$wb = $excel.ActiveWorkbook

# --- 1. Remove the "United States" sheet -----------------------------------
$wb.Worksheets.Item("United States").Delete() | Out-Null

# --- 2. Fix up the three remaining country sheets ---------------------------
$sheetNames = @("Australia", "Germany", "India")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Relabel the wafer / cell cost rows to their more precise names.
    $ws.Range("A4").Value = "Wafer (excl. polysilicon)"
    $ws.Range("A5").Value = "Cell Cost (excl. wafer)"

    # Relabel "Other material" (still sitting in row 12 at this point).
    $ws.Range("A12").Value = "Other material (e.g. front and back glass, encapsulant and others)"

    # Move the "Other material" row (row 12) up to directly follow
    # "Cell Cost" (i.e. become the new row 6), shifting the rows in
    # between (Overheads .. Labour, rows 6-11) down by one row.
    $block = $ws.Range("A6:D12").Value()

    $shifted = New-Object 'object[,]' 7,4
    for ($col = 0; $col -le 3; $col++) {
        $shifted[0, $col] = $block[7, $col + 1]
    }
    for ($row = 1; $row -le 6; $row++) {
        for ($col = 0; $col -le 3; $col++) {
            $shifted[$row, $col] = $block[$row, $col + 1]
        }
    }

    $ws.Range("A6:D12").Value = $shifted

    # The relocated label cell picks up the neighbouring numeric format
    # (as it would after an "Insert Cut Cells" move in real Excel).
    $ws.Range("A6").NumberFormat = "0.000"
}

# --- 3. Restore the saved selections ----------------------------------------
$wb.Worksheets.Item("Australia").Range("A13:XFD13").Select() | Out-Null
$wb.Worksheets.Item("Germany").Range("A13:XFD13").Select() | Out-Null
$wb.Worksheets.Item("India").Range("A5").Select() | Out-Null
$wb.Worksheets.Item("Australia").Activate() | Out-Null
